$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "cxq6hz_20170224_144343_ASIC_EEG"
$ws.Range("F1").Value = "hyy-调节6Hz_20170306_110203_ASIC_EEG"
$ws.Range("G1").Value = "hzj-调节6Hz_20170220_113105_ASIC_EEG"

$ws.Range("E2").Value = 0.5370919881305638
$ws.Range("F2").Value = 0.57575757575757569
$ws.Range("G2").Value = 0.53453453453453448

$ws.Range("E3").Value = 0.58309037900874627
$ws.Range("F3").Value = 0.60409556313993173
$ws.Range("G3").Value = 0.61612903225806459

$ws.Range("A1:G3").Select() | Out-Null
